$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header labels
$ws.Range("A1").Value = "Sector descripción"
$ws.Range("B1").Value = "Sector código"
$ws.Range("C1").Value = "Número Empresas"
$ws.Range("D1").Value = "Régimen"
$ws.Range("E1").Value = "Aragón'"
$ws.Range("F1").Value = "Dirección provincial nombre"
$ws.Range("G1").Value = "Mes y año"
$ws.Range("H1").Value = "Dirección provincial (código)"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:sector-descripcion"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "iaest-measure:numero-empresas"
$ws.Range("D2").Value = "iaest-measure:regimen"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("G2").Value = "iaest-measure:mes-y-ano"
$ws.Range("H2").Value = "null"

# Row 3 - medida/dim markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "null"

# Row 4 - xsd types / URI
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "xsd:string"
$ws.Range("H4").Value = "null"
